$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 9 new blank columns right after column A (i.e. at B:J), shifting
# all existing week-columns (old B:V) to the right by 9 (new K:AE).
$ws.Columns("B:J").Insert()

# ----- Row 1: header row with week-ending date labels (newest first) -----
$ws.Range("B1").Value2 = "Sep_08"
$ws.Range("C1").Value2 = "Aug_25"
$ws.Range("D1").Value2 = "Aug_04"
$ws.Range("E1").Value2 = "Jul_23"
$ws.Range("F1").Value2 = "Jul_17"
$ws.Range("G1").Value2 = "Jul_07"
$ws.Range("H1").Value2 = "Jun_30"
$ws.Range("I1").Value2 = "Jun_24"
$ws.Range("J1").Value2 = "Jun_16"

# ----- Fill the new (blank) data cells for every data row with "UN" -----
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
Write-Host "lastRow:" $lastRow
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Range("B$r`:J$r").Value2 = "UN"
}

# ----- Row 5 (Zacks Investment Research): two new rating changes -----
$ws.Range("B5").Value2 = "9/5/2019,Upgrades,Hold -> Buy,"
$ws.Range("B5").Interior.ColorIndex = 35
$ws.Range("C5").Value2 = "8/22/2019,Upgrades,Sell -> Hold,"
$ws.Range("C5").Interior.ColorIndex = 35

# ----- Row 11 (ValuEngine): one new rating change -----
$ws.Range("B11").Value2 = "9/6/2019,Upgrades,Hold -> Buy,"
$ws.Range("B11").Interior.ColorIndex = 35

# ----- Row 22 (BidaskClub): one new rating change -----
$ws.Range("B22").Value2 = "9/7/2019,Upgrades,Sell -> Hold,"
$ws.Range("B22").Interior.ColorIndex = 35

Write-Host "Done"
